$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (style) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-12 for the new columns I and J
$values = @(
    @(9, 9),
    @(7, 8),
    @(10, 10),
    @(9, 9),
    @(8, 8),
    @(6, 7),
    @(9, 9),
    @(1, 4),
    @(1, 3),
    @(6, 7),
    @(2, 3)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
